# Apply updated cryptos list values (prices & volume changes) to the sheet
# (GitHub Actions scheduled refresh of coinranking.com data)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds text-like values (e.g. "25.871.35") -- force text
# number format so Excel does not silently reinterpret them as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "25.871.35"
$ws.Cells.Item(2, 5).Value = "  -1.01%  "
$ws.Cells.Item(3, 4).Value = "1.638.34"
$ws.Cells.Item(3, 5).Value = "  -0.87%  "
$ws.Cells.Item(4, 4).Value = "0.9995"
$ws.Cells.Item(4, 5).Value = "  -1.34%  "
$ws.Cells.Item(5, 4).Value = "215.15"
$ws.Cells.Item(5, 5).Value = "  -0.25%  "
$ws.Cells.Item(6, 4).Value = "0.5039"
$ws.Cells.Item(6, 5).Value = "  -1.50%  "
$ws.Cells.Item(7, 4).Value = "1.001"
$ws.Cells.Item(7, 5).Value = "  -1.07%  "
$ws.Cells.Item(8, 2).Value = "Dogecoin"
$ws.Cells.Item(8, 3).Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Cells.Item(8, 4).Value = "0.06464"
$ws.Cells.Item(8, 5).Value = "  +0.44%  "
$ws.Cells.Item(9, 2).Value = "Cardano"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Cells.Item(9, 4).Value = "0.2572"
$ws.Cells.Item(9, 5).Value = "  -0.79%  "
$ws.Cells.Item(10, 4).Value = "19.52"
$ws.Cells.Item(10, 5).Value = "  -0.99%  "
$ws.Cells.Item(11, 4).Value = "0.07742"
$ws.Cells.Item(11, 5).Value = "  -0.62%  "
$ws.Cells.Item(12, 4).Value = "1.645.93"
$ws.Cells.Item(12, 5).Value = "  -0.40%  "
$ws.Cells.Item(13, 4).Value = "4.256"
$ws.Cells.Item(13, 5).Value = "  -0.49%  "
$ws.Cells.Item(14, 4).Value = "1.865.48"
$ws.Cells.Item(14, 5).Value = "  -0.84%  "
$ws.Cells.Item(15, 4).Value = "0.5458"
$ws.Cells.Item(15, 5).Value = "  -0.29%  "
$ws.Cells.Item(16, 4).Value = "0.0₅7917"
$ws.Cells.Item(16, 5).Value = "  -0.93%  "
$ws.Cells.Item(17, 4).Value = "64.24"
$ws.Cells.Item(17, 5).Value = "  +0.51%  "
$ws.Cells.Item(18, 4).Value = "25.896.59"
$ws.Cells.Item(18, 5).Value = "  -1.04%  "
$ws.Cells.Item(19, 4).Value = "1.001"
$ws.Cells.Item(19, 5).Value = "  -0.96%  "
$ws.Cells.Item(20, 4).Value = "202.10"
$ws.Cells.Item(20, 5).Value = "  -2.75%  "
$ws.Cells.Item(21, 4).Value = "4.378"
$ws.Cells.Item(21, 5).Value = "  -0.33%  "
$ws.Cells.Item(22, 4).Value = "9.889"
$ws.Cells.Item(22, 5).Value = "  -1.75%  "
$ws.Cells.Item(23, 4).Value = "5.979"
$ws.Cells.Item(23, 5).Value = "  -0.94%  "
$ws.Cells.Item(24, 5).Value = "  -0.91%  "
$ws.Cells.Item(25, 4).Value = "1.864"
$ws.Cells.Item(25, 5).Value = "  +0.49%  "
$ws.Cells.Item(26, 4).Value = "140.91"
$ws.Cells.Item(26, 5).Value = "  -2.28%  "
$ws.Cells.Item(27, 4).Value = "0.1135"
$ws.Cells.Item(27, 5).Value = "  -2.92%  "
$ws.Cells.Item(28, 5).Value = "  -0.93%  "
$ws.Cells.Item(29, 4).Value = "6.779"
$ws.Cells.Item(29, 5).Value = "  -2.36%  "
$ws.Cells.Item(30, 4).Value = "1.241"
$ws.Cells.Item(30, 5).Value = "  -0.11%  "
$ws.Cells.Item(31, 4).Value = "0.04927"
$ws.Cells.Item(31, 5).Value = "  -3.14%  "
$ws.Cells.Item(32, 4).Value = "3.271"
$ws.Cells.Item(32, 5).Value = "  -2.01%  "
$ws.Cells.Item(33, 4).Value = "3.197"
$ws.Cells.Item(33, 5).Value = "  -1.22%  "
$ws.Cells.Item(34, 4).Value = "1.545"
$ws.Cells.Item(34, 5).Value = "  -0.29%  "
$ws.Cells.Item(35, 4).Value = "2.361"
$ws.Cells.Item(35, 5).Value = "  +0.07%  "
$ws.Cells.Item(36, 2).Value = "ARBITRUM"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(36, 4).Value = "0.8913"
$ws.Cells.Item(36, 5).Value = "  -2.81%  "
$ws.Cells.Item(37, 2).Value = "MXToken"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(37, 4).Value = "2.620"
$ws.Cells.Item(37, 5).Value = "  -3.40%  "
$ws.Cells.Item(38, 4).Value = "1.145.51"
$ws.Cells.Item(38, 5).Value = "  -2.18%  "
$ws.Cells.Item(39, 4).Value = "0.5594"
$ws.Cells.Item(39, 5).Value = "  -1.74%  "
$ws.Cells.Item(40, 4).Value = "0.01563"
$ws.Cells.Item(40, 5).Value = "  -1.09%  "
$ws.Cells.Item(41, 4).Value = "1.002"
$ws.Cells.Item(41, 5).Value = "  -0.92%  "
$ws.Cells.Item(42, 5).Value = "  +0.39%  "
$ws.Cells.Item(43, 4).Value = "99.73"
$ws.Cells.Item(43, 5).Value = "  -0.53%  "
$ws.Cells.Item(44, 4).Value = "0.8052"
$ws.Cells.Item(44, 5).Value = "  -2.49%  "
$ws.Cells.Item(45, 4).Value = "1.776.85"
$ws.Cells.Item(45, 5).Value = "  -0.91%  "
$ws.Cells.Item(46, 5).Value = "  +3.05%  "
$ws.Cells.Item(47, 4).Value = "0.4508"
$ws.Cells.Item(47, 5).Value = "  -1.22%  "
$ws.Cells.Item(48, 4).Value = "1.008"
$ws.Cells.Item(48, 5).Value = "  -0.15%  "
$ws.Cells.Item(49, 4).Value = "54.66"
$ws.Cells.Item(49, 5).Value = "  -1.09%  "
$ws.Cells.Item(50, 4).Value = "0.05046"
$ws.Cells.Item(50, 5).Value = "  -0.80%  "
$ws.Cells.Item(51, 4).Value = "1.001"
$ws.Cells.Item(51, 5).Value = "  -0.89%  "
